# =====================================================================
# Update to drilling cycle variables sheet:
#  - add 'cycle_parting_off' variable block (1650-1658)
#  - add 'cycle_internal_threading' variable block (1700-1708)
#  - add 'cycle_OD_turning_chamfer_radius' variable block (1750+)
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-YellowFill($rng) {
    $rng.Interior.ColorIndex = 6
}
function Set-ThickBottomBorder($rng) {
    $bd = $rng.Borders.Item(9)
    $bd.LineStyle = 1
    $bd.Weight = -4138
    $bd.ColorIndex = 1
}

# --- Step 1: insert 25 blank rows before row 82 to make room for the new cycle blocks ---
$ws.Range("A82:C106").Insert(-4121) | Out-Null

# --- Step 2: populate values & labels ---
# Row 72
$ws.Range("A72").Value = 1614
$ws.Range("B72").Value = "rpm calculation"
$ws.Range("C72").Value = "cycle_ID_turning"

# Row 73
$ws.Range("A73").Value = 1650
$ws.Range("B73").Value = "Zstart"
$ws.Range("C73").Value = "cycle_parting_off"

# Row 74
$ws.Range("A74").Value = 1651
$ws.Range("B74").Value = "diameter A"
$ws.Range("C74").Value = "cycle_parting_off"

# Row 75
$ws.Range("A75").Value = 1652
$ws.Range("B75").Value = "diameter B"
$ws.Range("C75").Value = "cycle_parting_off"

# Row 76
$ws.Range("A76").Value = 1653
$ws.Range("B76").Value = "tool width"
$ws.Range("C76").Value = "cycle_parting_off"

# Row 77
$ws.Range("A77").Value = 1654
$ws.Range("B77").Value = "feed per rev"
$ws.Range("C77").Value = "cycle_parting_off"

# Row 78
$ws.Range("A78").Value = 1655
$ws.Range("B78").Value = "spindle speed"
$ws.Range("C78").Value = "cycle_parting_off"

# Row 79
$ws.Range("A79").Value = 1656
$ws.Range("B79").Value = "pecking depth"
$ws.Range("C79").Value = "cycle_parting_off"

# Row 80
$ws.Range("A80").Value = 1657
$ws.Range("B80").Value = "retract amount"
$ws.Range("C80").Value = "cycle_parting_off"

# Row 81
$ws.Range("A81").Value = 1658
$ws.Range("B81").Value = "dwell time"
$ws.Range("C81").Value = "cycle_parting_off"

# Row 82
$ws.Range("A82").Value = 1700
$ws.Range("B82").Value = "Z1"
$ws.Range("C82").Value = "cycle_internal_threading"

# Row 83
$ws.Range("A83").Value = 1701
$ws.Range("B83").Value = "Z2"
$ws.Range("C83").Value = "cycle_internal_threading"

# Row 84
$ws.Range("A84").Value = 1702
$ws.Range("B84").Value = "diameter A"
$ws.Range("C84").Value = "cycle_internal_threading"

# Row 85
$ws.Range("A85").Value = 1703
$ws.Range("B85").Value = "diameter B"
$ws.Range("C85").Value = "cycle_internal_threading"

# Row 86
$ws.Range("A86").Value = 1704
$ws.Range("B86").Value = "pitch"
$ws.Range("C86").Value = "cycle_internal_threading"

# Row 87
$ws.Range("A87").Value = 1705
$ws.Range("B87").Value = "depth per pass"
$ws.Range("C87").Value = "cycle_internal_threading"

# Row 88
$ws.Range("A88").Value = 1706
$ws.Range("B88").Value = "spindle speed rev/min"
$ws.Range("C88").Value = "cycle_internal_threading"

# Row 89
$ws.Range("A89").Value = 1707
$ws.Range("B89").Value = "full thread depth beyond thread peak"
$ws.Range("C89").Value = "cycle_internal_threading"

# Row 90
$ws.Range("A90").Value = 1708
$ws.Range("B90").Value = "Z clearance"
$ws.Range("C90").Value = "cycle_internal_threading"

# Row 91
$ws.Range("A91").Value = 1750
$ws.Range("C91").Value = "cycle_OD_turning_chamfer_radius"

# Row 92
$ws.Range("C92").Value = "cycle_OD_turning_chamfer_radius"

# Row 93
$ws.Range("C93").Value = "cycle_OD_turning_chamfer_radius"

# Row 94
$ws.Range("C94").Value = "cycle_OD_turning_chamfer_radius"

# Row 95
$ws.Range("C95").Value = "cycle_OD_turning_chamfer_radius"

# Row 96
$ws.Range("C96").Value = "cycle_OD_turning_chamfer_radius"

# Row 97
$ws.Range("C97").Value = "cycle_OD_turning_chamfer_radius"

# Row 98
$ws.Range("C98").Value = "cycle_OD_turning_chamfer_radius"

# Row 99

# Row 100

# Row 101

# Row 102

# Row 103

# Row 104

# Row 105

# Row 106

# --- Step 3: apply cell shading / borders to match the cycle-block formatting ---
Set-YellowFill($ws.Range("A72"))
Set-ThickBottomBorder($ws.Range("A72"))
Set-ThickBottomBorder($ws.Range("B72"))
Set-ThickBottomBorder($ws.Range("C72"))

Set-YellowFill($ws.Range("A73"))

Set-YellowFill($ws.Range("A74"))

Set-YellowFill($ws.Range("A75"))

Set-YellowFill($ws.Range("A76"))

Set-YellowFill($ws.Range("A77"))

Set-YellowFill($ws.Range("A78"))

Set-YellowFill($ws.Range("A79"))

Set-YellowFill($ws.Range("A80"))

Set-YellowFill($ws.Range("A81"))
Set-ThickBottomBorder($ws.Range("A81"))
Set-ThickBottomBorder($ws.Range("B81"))
Set-ThickBottomBorder($ws.Range("C81"))

Set-YellowFill($ws.Range("A82"))

Set-YellowFill($ws.Range("A83"))

Set-YellowFill($ws.Range("A84"))

Set-YellowFill($ws.Range("A85"))

Set-YellowFill($ws.Range("A86"))

Set-YellowFill($ws.Range("A87"))

Set-YellowFill($ws.Range("A88"))

Set-YellowFill($ws.Range("A89"))

Set-YellowFill($ws.Range("A90"))
Set-ThickBottomBorder($ws.Range("A90"))
Set-ThickBottomBorder($ws.Range("B90"))
Set-ThickBottomBorder($ws.Range("C90"))

Set-YellowFill($ws.Range("A91"))

Set-YellowFill($ws.Range("A92"))

Set-YellowFill($ws.Range("A93"))

Set-YellowFill($ws.Range("A94"))

Set-YellowFill($ws.Range("A95"))

Set-YellowFill($ws.Range("A96"))

Set-YellowFill($ws.Range("A97"))

Set-YellowFill($ws.Range("A98"))

Set-YellowFill($ws.Range("A99"))

Set-YellowFill($ws.Range("A100"))

Set-YellowFill($ws.Range("A101"))

Set-YellowFill($ws.Range("A102"))

Set-YellowFill($ws.Range("A103"))

Set-YellowFill($ws.Range("A104"))

Set-YellowFill($ws.Range("A105"))

Set-YellowFill($ws.Range("A106"))

# --- Step 4: restore view selection (active cell) ---
$ws.Range("E104").Select() | Out-Null

Write-Host "edit complete"
